$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Move the "_GoBack" bookmark from the end of the fourth paragraph
#    ("Fjärde paragrafen...") up into the title ("Referensdokument"),
#    splitting "Referensdokument" into "Referensd" | "okument" with
#    the bookmark sitting between the two runs.
# ------------------------------------------------------------------
$titleRange = $d.Content
$null = $titleRange.Find.Execute("Referensdokument")
$splitOffset = $titleRange.Start + "Referensd".Length

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$splitPoint = $d.Range($splitOffset, $splitOffset)
$d.Bookmarks.Add("_GoBack", $splitPoint)

# ------------------------------------------------------------------
# 2) "no indrag on title" -- stop the Title ("Rubrik") style from
#    inheriting Normal's first-line indent.
# ------------------------------------------------------------------
$titleStyle = $d.Styles("Title")
$titleStyle.ParagraphFormat.FirstLineIndent = 0
